# Updates cryptos list values/percentages (and a few reordered rows) per the
# upstream data refresh. Column D ("Price") cells are forced to Text via a
# temporary "@" number format so numeric-looking strings (e.g. "237.75")
# are not silently converted to numbers by Excel's type inference; the
# format is cleared right after so no stray cell style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "94.291.38"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +2.54%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.498.07"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +6.01%  "
$ws.Range("E4").Value = "  -0.09%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "237.75"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +4.71%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "627.82"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +2.80%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "1.44"
$cell.ClearFormats()
$ws.Range("E7").Value = "  +6.89%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.396"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +5.88%  "
$ws.Range("E9").Value = "  -0.14%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "1.01"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +11.47%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "3.496.88"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +5.94%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "43.31"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +5.44%  "
$ws.Range("E13").Value = "  +6.74%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.27"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +6.94%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.155.38"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +6.08%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "94.081.30"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("E17").Value = "  +5.42%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "8.38"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +7.59%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "3.498.39"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +6.02%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "12.75"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +19.02%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "18.14"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +8.69%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.504"
$cell.ClearFormats()
$ws.Range("E22").Value = "  +16.57%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "518.91"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +8.24%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "3.37"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +5.20%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "6.76"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +12.87%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.0000186"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +4.62%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "96.06"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +8.24%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "12.24"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +7.87%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "3.679.20"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +5.61%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.93"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +13.04%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "11.55"
$cell.ClearFormats()
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("E32").Value = "  +0.03%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.140"
$cell.ClearFormats()
$ws.Range("E33").Value = "  +5.69%  "
$ws.Range("B34").Value = "Cronos"
$ws.Range("C34").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.181"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +7.72%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.993"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.563"
$cell.ClearFormats()
$ws.Range("E36").Value = "  +9.36%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "29.78"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +7.27%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "577.80"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +12.93%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.46"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +10.46%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "7.58"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +5.88%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.929"
$cell.ClearFormats()
$ws.Range("E42").Value = "  +8.18%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.149"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +4.35%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0426"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +9.54%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "23.75"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -0.82%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.71"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +5.31%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "5.54"
$cell.ClearFormats()
$ws.Range("E47").Value = "  +6.14%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.57"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +0.60%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.ClearFormats()
$ws.Range("E49").Value = "  +5.70%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "53.55"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +3.75%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "8.22"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +5.38%  "
